# Daily attendance processing - 2026-01-20 07:43:33
# Swap the order of the two comma-separated "Recorded By" names in column G
# for the specific combinations that need reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($text -eq "System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, System"
    }
    elseif ($text -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
}
